$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add placeholder label for DOI of dataset in A3
$ws.Range("A3").Value = "[DOI of dataset]"

# Move the "20" value up in column B: B4 becomes blank, B5 becomes 20
$ws.Range("B4").ClearContents()
$ws.Range("B5").Value = 20

# Update the last active selection to C4
$ws.Range("C4").Select()
